$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.531.51'
$ws.Range('E2').Value = '  -1.16%  '
$ws.Range('D3').Value = '2.995.64'
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''600.48'
$ws.Range('E5').Value = '  +3.30%  '
$ws.Range('D6').Value = '''144.40'
$ws.Range('E6').Value = '  -1.55%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '''0.520'
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').Value = '2.992.81'
$ws.Range('E9').Value = '  -0.32%  '
$ws.Range('E10').Value = '  -0.71%  '
$ws.Range('D11').Value = '''6.05'
$ws.Range('E11').Value = '  +7.34%  '
$ws.Range('D12').Value = '''0.456'
$ws.Range('E12').Value = '  +3.85%  '
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('D14').Value = '''34.40'
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('E15').Value = '  +2.37%  '
$ws.Range('D16').Value = '3.492.81'
$ws.Range('E16').Value = '  -0.40%  '
$ws.Range('D17').Value = '''6.97'
$ws.Range('E17').Value = '  -0.70%  '
$ws.Range('D18').Value = '61.476.52'
$ws.Range('E18').Value = '  -1.30%  '
$ws.Range('D19').Value = '2.996.81'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').Value = '''451.85'
$ws.Range('E20').Value = '  -0.75%  '
$ws.Range('D21').Value = '''14.01'
$ws.Range('E21').Value = '  +1.15%  '
$ws.Range('E22').Value = '  +1.41%  '
$ws.Range('D23').Value = '''7.33'
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('D24').Value = '''81.58'
$ws.Range('E24').Value = '  +2.17%  '
$ws.Range('D25').Value = '''10.82'
$ws.Range('E25').Value = '  +7.34%  '
$ws.Range('E26').Value = '  -2.84%  '
$ws.Range('D27').Value = '''11.98'
$ws.Range('E27').Value = '  -2.41%  '
$ws.Range('E28').Value = '  +0.24%  '
$ws.Range('E29').Value = '  +3.10%  '
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('D31').Value = '''7.20'
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('D32').Value = '''2.06'
$ws.Range('E32').Value = '  -1.29%  '
$ws.Range('D33').Value = '''27.35'
$ws.Range('E33').Value = '  +1.79%  '
$ws.Range('E34').Value = '  +3.46%  '
$ws.Range('D35').Value = '0.0₃0830'
$ws.Range('E35').Value = '  +5.28%  '
$ws.Range('E36').Value = '  -0.42%  '
$ws.Range('E37').Value = '  +1.36%  '
$ws.Range('D38').Value = '''9.21'
$ws.Range('E38').Value = '  +2.64%  '
$ws.Range('D39').Value = '''50.41'
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('E41').Value = '  +11.18%  '
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').Value = '''397.64'
$ws.Range('E43').Value = '  -4.24%  '
$ws.Range('D44').Value = '''39.77'
$ws.Range('E44').Value = '  +4.50%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '''0.0354'
$ws.Range('E45').Value = '  +0.20%  '
$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').Value = '''0.271'
$ws.Range('E46').Value = '  -1.43%  '
$ws.Range('D47').Value = '2.692.98'
$ws.Range('E47').Value = '  -2.66%  '
$ws.Range('D48').Value = '''130.82'
$ws.Range('E48').Value = '  +1.89%  '
$ws.Range('E50').Value = '  -0.22%  '
$ws.Range('E51').Value = '  +2.16%  '
